# Applies the "05-11-2023 14:45" script update to the Serbia Prva Liga
# 2023-2024 betting-odds sheet:
#   1) A handful of existing match rows had their F:V (match/odds/url)
#      payload reshuffled among same-kick-off-time rows (A:E — index,
#      country, tournament, season, date — stay put).
#   2) Three brand-new matches are appended as rows 103-105.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($row, $startCol, $endCol) {
    $vals = @()
    for ($c = $startCol; $c -le $endCol; $c++) {
        $vals += $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($row, $startCol, $vals) {
    $c = $startCol
    foreach ($v in $vals) {
        $ws.Cells.Item($row, $c).Value2 = $v
        $c = $c + 1
    }
}

# --- Step 1: re-shuffle F:V (columns 6..22) across the affected rows ---
# Capture every "source" row's current payload BEFORE any writes happen,
# since several groups are 3-way rotations (writes would otherwise
# clobber a value still needed for a later read).
$src2  = Get-RowValues 2  6 22
$src4  = Get-RowValues 4  6 22

$src25 = Get-RowValues 25 6 22
$src26 = Get-RowValues 26 6 22

$src30 = Get-RowValues 30 6 22
$src31 = Get-RowValues 31 6 22

$src42 = Get-RowValues 42 6 22
$src43 = Get-RowValues 43 6 22
$src44 = Get-RowValues 44 6 22

$src49 = Get-RowValues 49 6 22
$src50 = Get-RowValues 50 6 22

$src52 = Get-RowValues 52 6 22
$src53 = Get-RowValues 53 6 22
$src54 = Get-RowValues 54 6 22

$src79 = Get-RowValues 79 6 22
$src80 = Get-RowValues 80 6 22

# Now apply the target permutation for each group.
Set-RowValues 2  6 $src4
Set-RowValues 4  6 $src2

Set-RowValues 25 6 $src26
Set-RowValues 26 6 $src25

Set-RowValues 30 6 $src31
Set-RowValues 31 6 $src30

Set-RowValues 42 6 $src43
Set-RowValues 43 6 $src44
Set-RowValues 44 6 $src42

Set-RowValues 49 6 $src50
Set-RowValues 50 6 $src49

Set-RowValues 52 6 $src54
Set-RowValues 53 6 $src52
Set-RowValues 54 6 $src53

Set-RowValues 79 6 $src80
Set-RowValues 80 6 $src79

# --- Step 2: append the 3 new matches as rows 103-105 ---
# Clone row 102's formatting (bold/border index cell in A, date format in
# E) onto the new rows before filling in values.
$ws.Range("A102:V102").Copy()
$ws.Range("A103:V103").PasteSpecial(-4122)
$ws.Range("A104:V104").PasteSpecial(-4122)
$ws.Range("A105:V105").PasteSpecial(-4122)

$newRows = @(
    @(102, "serbia", "prva-liga", "2023-2024", 45235.54166666666, `
      "Graficar Beograd", 5, "Sloboda", 1, `
      1.74, "04/11/2023 01:13", 1.6, "05/11/2023 12:45", `
      3.15, "04/11/2023 01:13", 3.49, "05/11/2023 12:45", `
      3.92, "04/11/2023 01:13", 5.14, "05/11/2023 12:45", `
      "https://www.betexplorer.com/football/serbia/prva-liga/graficar-beograd-sloboda/rRCAlpTT/"),
    @(103, "serbia", "prva-liga", "2023-2024", 45235.54166666666, `
      "OFK Beograd", 2, "Macva", 1, `
      1.28, "04/11/2023 01:13", 1.38, "05/11/2023 12:42", `
      4.14, "04/11/2023 01:13", 4.28, "05/11/2023 12:42", `
      7.83, "04/11/2023 01:13", 6.82, "05/11/2023 12:42", `
      "https://www.betexplorer.com/football/serbia/prva-liga/ofk-beograd-macva-sabac/Kn4xp4yo/"),
    @(104, "serbia", "prva-liga", "2023-2024", 45235.54166666666, `
      "Tekstilac Odzaci", 1, "Vrsac", 0, `
      1.74, "04/11/2023 01:13", 1.88, "05/11/2023 12:53", `
      2.92, "04/11/2023 01:13", 2.86, "05/11/2023 12:54", `
      4.32, "04/11/2023 01:13", 4.5, "05/11/2023 12:53", `
      "https://www.betexplorer.com/football/serbia/prva-liga/tekstilac-odzaci-vrsac/YkhK8mEo/")
)

$r = 103
foreach ($row in $newRows) {
    Set-RowValues $r 1 $row
    $r = $r + 1
}
